# Add the "2021年" (2021) data row to the bottom of the table, one row
# below the existing last row (row 11, "2020年").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carry the row-11 formatting (bold/centered/bordered year style in column
# A, plain style in the data columns) down onto the new row 12 first, so
# the new cells line up with the rest of the table.
$ws.Range("A11:BK11").Copy() | Out-Null
$ws.Range("A12:BK12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row label
$ws.Range("A12").Value = "2021年"

# Data values for 2021 (columns without a value in the source data are
# left blank, matching the sparse columns already present in the table).
$ws.Range("C12").Value = 233
$ws.Range("D12").Value = 145
$ws.Range("H12").Value = 995
$ws.Range("I12").Value = 68
$ws.Range("J12").Value = 2214
$ws.Range("K12").Value = 19
$ws.Range("L12").Value = 2034
$ws.Range("M12").Value = 295
$ws.Range("N12").Value = 68
$ws.Range("Q12").Value = 2
$ws.Range("R12").Value = 4
$ws.Range("S12").Value = 31
$ws.Range("U12").Value = 21
$ws.Range("V12").Value = 16
$ws.Range("X12").Value = 361
$ws.Range("Y12").Value = 498
$ws.Range("Z12").Value = 368
$ws.Range("AA12").Value = 460
$ws.Range("AC12").Value = 21
$ws.Range("AD12").Value = 38
$ws.Range("AE12").Value = 22
$ws.Range("AF12").Value = 263
$ws.Range("AG12").Value = 1567
$ws.Range("AH12").Value = 68
$ws.Range("AI12").Value = 6
$ws.Range("AJ12").Value = 2
$ws.Range("AK12").Value = 23
$ws.Range("AN12").Value = 70
$ws.Range("AO12").Value = 29
$ws.Range("AP12").Value = 90
$ws.Range("AR12").Value = 13
$ws.Range("AS12").Value = 169
$ws.Range("AV12").Value = 5
$ws.Range("AW12").Value = 272
$ws.Range("AX12").Value = 67
$ws.Range("AY12").Value = 279
$ws.Range("AZ12").Value = 8
$ws.Range("BB12").Value = 11
$ws.Range("BC12").Value = 176
$ws.Range("BD12").Value = 26
$ws.Range("BE12").Value = 484
$ws.Range("BG12").Value = 11617
$ws.Range("BH12").Value = 17
$ws.Range("BJ12").Value = 36
$ws.Range("BK12").Value = 23
